$wb = $excel.ActiveWorkbook

# --- Update the transport_solution lookup sheet ---
$ts = $wb.Worksheets.Item("transport_solution")

$ts.Range("A1").Value = "PBS"
$ts.Range("B1").Value = "http://purl.obolibrary.org/obo/OBI_0100046"

$ts.Range("A2").Value = "Saline (Buffered)"
$ts.Range("B2").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000154"

$ts.Range("A3").Value = "UWS"
$ts.Range("B3").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000151"

$ts.Range("A4").Value = "DMEM"
$ts.Range("B4").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185409"

$ts.Range("A5").Value = "Miltenyi Tissue Preservation Buffer"
$ts.Range("B5").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000150"

$ts.Range("A6").Value = "NBF (Neutral Buffered Formalin)"
$ts.Range("B6").Value = "http://purl.obolibrary.org/obo/OBIB_0000213"

$ts.Range("A7").Value = "Unknown"
$ts.Range("B7").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998"

# Rows 8-11 are unchanged (RPMI, None, HTK, Belzer MPS/KPS)

# --- Update the .metadata sheet's pav:createdOn value ---
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2023-08-04T07:36:03-07:00"
